$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 990.4805
$ws.Range("I129").Value = 951.8
$ws.Range("J129").Value = 996.2537
$ws.Range("K129").Value = 2855.4
$ws.Range("L129").Value = 2988.7611
$ws.Range("M129").Value = 2144.6
$ws.Range("N129").Value = -12988.7611

$ws.Range("H137").Value = 1681.5927
$ws.Range("I137").Value = 1455.5555
$ws.Range("J137").Value = 2133.6667
$ws.Range("K137").Value = 4366.666499999999
$ws.Range("L137").Value = 6401.000100000001
$ws.Range("M137").Value = -1816.666499999999
$ws.Range("N137").Value = -11501.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15875744
$ws.Range("I32").Value = 2429.2075
$ws.Range("J32").Value = 100004310
$ws.Range("K32").Value = 2429.2075
$ws.Range("L32").Value = 100004310
$ws.Range("M32").Value = -2142.2075
$ws.Range("N32").Value = -100004884

$ws.Range("H61").Value = 25001504
$ws.Range("I61").Value = 25001504
$ws.Range("K61").Value = 25001504
$ws.Range("M61").Value = -25001292

$ws.Range("H132").Value = 2101993.2
$ws.Range("I132").Value = 1128.5
$ws.Range("J132").Value = 14707182
$ws.Range("K132").Value = 3385.5
$ws.Range("L132").Value = 44121546
$ws.Range("M132").Value = -855.5
$ws.Range("N132").Value = -44126606

$ws.Range("H136").Value = 25001504
$ws.Range("I136").Value = 25001504
$ws.Range("K136").Value = 75004512
$ws.Range("M136").Value = -75001962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1754.0625
$ws.Range("I99").Value = 1270
$ws.Range("K99").Value = 1270
$ws.Range("M99").Value = 228

$ws.Range("H132").Value = 48853.332
$ws.Range("J132").Value = 48853.332
$ws.Range("L132").Value = 48853.332
$ws.Range("N132").Value = -58973.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187

$ws.Range("H6").Value = 125472.625
$ws.Range("I6").Value = 333517
$ws.Range("J6").Value = 646
$ws.Range("K6").Value = 333517
$ws.Range("L6").Value = 646
$ws.Range("M6").Value = -333404
$ws.Range("N6").Value = -872

$ws.Range("H12").Value = 629.3333
$ws.Range("I12").Value = 265.66666
$ws.Range("K12").Value = 265.66666
$ws.Range("M12").Value = -95.66665999999998

$ws.Range("H31").Value = 2223237
$ws.Range("I31").Value = 2416318.5
$ws.Range("J31").Value = 2799.5
$ws.Range("K31").Value = 2416318.5
$ws.Range("L31").Value = 2799.5
$ws.Range("M31").Value = -2416023.5
$ws.Range("N31").Value = -3389.5

$ws.Range("H34").Value = 2223237
$ws.Range("I34").Value = 2416318.5
$ws.Range("J34").Value = 2799.5
$ws.Range("K34").Value = 2416318.5
$ws.Range("L34").Value = 2799.5
$ws.Range("M34").Value = -2416116.5
$ws.Range("N34").Value = -3203.5

$ws.Range("H134").Value = 1894.2
$ws.Range("I134").Value = 1650
$ws.Range("J134").Value = 2260.5
$ws.Range("K134").Value = 4950
$ws.Range("L134").Value = 6781.5
$ws.Range("M134").Value = -2415
$ws.Range("N134").Value = -11851.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 43214100
$ws.Range("I5").Value = 53030656
$ws.Range("K5").Value = 159091968
$ws.Range("M5").Value = -159091856

$ws.Range("H106").Value = 18666.666
$ws.Range("J106").Value = 18666.666
$ws.Range("L106").Value = 55999.99800000001
$ws.Range("N106").Value = -57891.99800000001

$ws.Range("H122").Value = 13025525
$ws.Range("J122").Value = 6176.3613
$ws.Range("L122").Value = 55587.25169999999
$ws.Range("N122").Value = -60487.25169999999

$ws.Range("H131").Value = 817.27
$ws.Range("I131").Value = 657.5
$ws.Range("J131").Value = 823.92706
$ws.Range("K131").Value = 1972.5
$ws.Range("L131").Value = 2471.78118
$ws.Range("M131").Value = 3067.5
$ws.Range("N131").Value = -12551.78118

$ws.Range("H132").Value = 5135.269
$ws.Range("I132").Value = 738.4167
$ws.Range("J132").Value = 8904
$ws.Range("K132").Value = 6645.7503
$ws.Range("L132").Value = 80136
$ws.Range("M132").Value = -4115.7503
$ws.Range("N132").Value = -85196

$ws.Range("H135").Value = 43214100
$ws.Range("I135").Value = 53030656
$ws.Range("K135").Value = 477275904
$ws.Range("M135").Value = -477273369

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 6666
$ws.Range("J62").Value = 6666
$ws.Range("L62").Value = 6666
$ws.Range("N62").Value = -8038

$ws.Range("H65").Value = 6666
$ws.Range("J65").Value = 6666
$ws.Range("L65").Value = 19998
$ws.Range("N65").Value = -26862

$ws.Range("H132").Value = 11035.546
$ws.Range("I132").Value = 2243.4443
$ws.Range("J132").Value = 50600
$ws.Range("K132").Value = 6730.3329
$ws.Range("L132").Value = 151800
$ws.Range("M132").Value = -4200.3329
$ws.Range("N132").Value = -156860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2037.6666
$ws.Range("I7").Value = 1207.5714
$ws.Range("J7").Value = 3199.8
$ws.Range("K7").Value = 1207.5714
$ws.Range("L7").Value = 3199.8
$ws.Range("M7").Value = -1095.5714
$ws.Range("N7").Value = -3423.8

$ws.Range("H63").Value = 12324.333
$ws.Range("J63").Value = 12324.333
$ws.Range("L63").Value = 12324.333
$ws.Range("N63").Value = -13822.333

$ws.Range("H66").Value = 12324.333
$ws.Range("J66").Value = 12324.333
$ws.Range("L66").Value = 36972.999
$ws.Range("N66").Value = -44460.999

$ws.Range("H68").Value = 1359.8334
$ws.Range("I68").Value = 1391.9286
$ws.Range("J68").Value = 1247.5
$ws.Range("K68").Value = 1391.9286
$ws.Range("L68").Value = 1247.5
$ws.Range("M68").Value = -642.9286
$ws.Range("N68").Value = -2745.5

$ws.Range("H71").Value = 1359.8334
$ws.Range("I71").Value = 1391.9286
$ws.Range("J71").Value = 1247.5
$ws.Range("K71").Value = 6959.643
$ws.Range("L71").Value = 6237.5
$ws.Range("M71").Value = -3215.643
$ws.Range("N71").Value = -13725.5

$ws.Range("H93").Value = 1196.875
$ws.Range("I93").Value = 961.6667
$ws.Range("J93").Value = 1588.8889
$ws.Range("K93").Value = 961.6667
$ws.Range("L93").Value = 1588.8889
$ws.Range("M93").Value = 286.3333
$ws.Range("N93").Value = -4084.8889

$ws.Range("H126").Value = 2037.6666
$ws.Range("I126").Value = 1207.5714
$ws.Range("J126").Value = 3199.8
$ws.Range("K126").Value = 3622.7142
$ws.Range("L126").Value = 9599.400000000001
$ws.Range("M126").Value = -1152.7142
$ws.Range("N126").Value = -14539.4

$ws.Range("H132").Value = 43966304
$ws.Range("I132").Value = 81633990
$ws.Range("J132").Value = 20666.166
$ws.Range("K132").Value = 244901970
$ws.Range("L132").Value = 61998.49800000001
$ws.Range("M132").Value = -244899440
$ws.Range("N132").Value = -67058.49800000001

$ws.Range("H136").Value = 79367000
$ws.Range("I136").Value = 60152716
$ws.Range("J136").Value = 125000930
$ws.Range("K136").Value = 180458148
$ws.Range("L136").Value = 375002790
$ws.Range("M136").Value = -180455598
$ws.Range("N136").Value = -375007890
